$wb = $excel.ActiveWorkbook

# --- ALC row 15 (Leve Item ID 44146) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 649.5821
$ws.Cells.Item(15, 9).Value = 649.5821
$ws.Cells.Item(15, 11).Value = 1948.7463
$ws.Cells.Item(15, 13).Value = -1779.7463

# --- ALC row 33 (Leve Item ID 5512) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 151.36363
$ws.Cells.Item(33, 9).Value = 145.8
$ws.Cells.Item(33, 10).Value = 207
$ws.Cells.Item(33, 11).Value = 145.8
$ws.Cells.Item(33, 12).Value = 207
$ws.Cells.Item(33, 13).Value = 83.19999999999999
$ws.Cells.Item(33, 14).Value = -665

# --- ALC row 61 (Leve Item ID 4604) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(61, 8).Value = 16000
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 13).Value = ""

# --- ALC row 74 (Leve Item ID 5507) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 4978.6
$ws.Cells.Item(74, 9).Value = 4754.857
$ws.Cells.Item(74, 10).Value = 5500.6665
$ws.Cells.Item(74, 11).Value = 4754.857
$ws.Cells.Item(74, 12).Value = 5500.6665
$ws.Cells.Item(74, 13).Value = -3818.857
$ws.Cells.Item(74, 14).Value = -7372.6665

# --- ALC row 77 (Leve Item ID 5507) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 4978.6
$ws.Cells.Item(77, 9).Value = 4754.857
$ws.Cells.Item(77, 10).Value = 5500.6665
$ws.Cells.Item(77, 11).Value = 23774.285
$ws.Cells.Item(77, 12).Value = 27503.3325
$ws.Cells.Item(77, 13).Value = -19094.285
$ws.Cells.Item(77, 14).Value = -36863.3325

# --- ALC row 80 (Leve Item ID 12605) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 10537.333
$ws.Cells.Item(80, 9).Value = 637.5
$ws.Cells.Item(80, 10).Value = 15487.25
$ws.Cells.Item(80, 11).Value = 1912.5
$ws.Cells.Item(80, 12).Value = 46461.75
$ws.Cells.Item(80, 13).Value = -914.5
$ws.Cells.Item(80, 14).Value = -48457.75

# --- ALC row 83 (Leve Item ID 12605) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(83, 8).Value = 10537.333
$ws.Cells.Item(83, 9).Value = 637.5
$ws.Cells.Item(83, 10).Value = 15487.25
$ws.Cells.Item(83, 11).Value = 5737.5
$ws.Cells.Item(83, 12).Value = 139385.25
$ws.Cells.Item(83, 13).Value = -745.5
$ws.Cells.Item(83, 14).Value = -149369.25

# --- ALC row 137 (Leve Item ID 44013) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 559705.75
$ws.Cells.Item(137, 9).Value = 1567.75
$ws.Cells.Item(137, 10).Value = 2420165.8
$ws.Cells.Item(137, 11).Value = 4703.25
$ws.Cells.Item(137, 12).Value = 7260497.399999999
$ws.Cells.Item(137, 13).Value = -2153.25
$ws.Cells.Item(137, 14).Value = -7265597.399999999

# --- ARM row 32 (Leve Item ID 44147) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9161.022999999999
$ws.Cells.Item(32, 9).Value = 3310.138
$ws.Cells.Item(32, 11).Value = 3310.138
$ws.Cells.Item(32, 13).Value = -3023.138

# --- ARM row 45 (Leve Item ID 27714) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 8931481
$ws.Cells.Item(45, 9).Value = 2405
$ws.Cells.Item(45, 10).Value = 15628288
$ws.Cells.Item(45, 11).Value = 2405
$ws.Cells.Item(45, 12).Value = 15628288
$ws.Cells.Item(45, 13).Value = -2028
$ws.Cells.Item(45, 14).Value = -15629042

# --- ARM row 74 (Leve Item ID 44000) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 40742.96
$ws.Cells.Item(74, 9).Value = 49181.76
$ws.Cells.Item(74, 11).Value = 49181.76
$ws.Cells.Item(74, 13).Value = -48307.76

# --- ARM row 77 (Leve Item ID 44000) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 40742.96
$ws.Cells.Item(77, 9).Value = 49181.76
$ws.Cells.Item(77, 11).Value = 245908.8
$ws.Cells.Item(77, 13).Value = -241540.8

# --- ARM row 80 (Leve Item ID 10667) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(80, 8).Value = 77900
$ws.Cells.Item(80, 10).Value = 77900
$ws.Cells.Item(80, 12).Value = 77900
$ws.Cells.Item(80, 14).Value = -79896

# --- ARM row 83 (Leve Item ID 10667) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(83, 8).Value = 77900
$ws.Cells.Item(83, 10).Value = 77900
$ws.Cells.Item(83, 12).Value = 233700
$ws.Cells.Item(83, 14).Value = -243684

# --- ARM row 107 (Leve Item ID 25645) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(107, 8).Value = 62159
$ws.Cells.Item(107, 10).Value = 62159
$ws.Cells.Item(107, 12).Value = 62159
$ws.Cells.Item(107, 14).Value = -69839

# --- BSM row 86 (Leve Item ID 12526) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3373.68
$ws.Cells.Item(86, 9).Value = 3470.3333
$ws.Cells.Item(86, 10).Value = 3228.7
$ws.Cells.Item(86, 11).Value = 3470.3333
$ws.Cells.Item(86, 12).Value = 3228.7
$ws.Cells.Item(86, 13).Value = -2347.3333
$ws.Cells.Item(86, 14).Value = -5474.7

# --- BSM row 89 (Leve Item ID 12526) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 3373.68
$ws.Cells.Item(89, 9).Value = 3470.3333
$ws.Cells.Item(89, 10).Value = 3228.7
$ws.Cells.Item(89, 11).Value = 17351.6665
$ws.Cells.Item(89, 12).Value = 16143.5
$ws.Cells.Item(89, 13).Value = -11735.6665
$ws.Cells.Item(89, 14).Value = -27375.5

# --- BSM row 107 (Leve Item ID 27706) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 5558550.5
$ws.Cells.Item(107, 9).Value = 7695315
$ws.Cells.Item(107, 10).Value = 2962.4
$ws.Cells.Item(107, 11).Value = 7695315
$ws.Cells.Item(107, 12).Value = 2962.4
$ws.Cells.Item(107, 13).Value = -7693395
$ws.Cells.Item(107, 14).Value = -6802.4

# --- CRP row 31 (Leve Item ID 44023) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3158.353
$ws.Cells.Item(31, 9).Value = 1657
$ws.Cells.Item(31, 10).Value = 5303.143
$ws.Cells.Item(31, 11).Value = 1657
$ws.Cells.Item(31, 12).Value = 5303.143
$ws.Cells.Item(31, 13).Value = -1362
$ws.Cells.Item(31, 14).Value = -5893.143

# --- CRP row 34 (Leve Item ID 44023) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3158.353
$ws.Cells.Item(34, 9).Value = 1657
$ws.Cells.Item(34, 10).Value = 5303.143
$ws.Cells.Item(34, 11).Value = 1657
$ws.Cells.Item(34, 12).Value = 5303.143
$ws.Cells.Item(34, 13).Value = -1455
$ws.Cells.Item(34, 14).Value = -5707.143

# --- CRP row 86 (Leve Item ID 12584) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 17860142
$ws.Cells.Item(86, 10).Value = 6000
$ws.Cells.Item(86, 12).Value = 6000
$ws.Cells.Item(86, 14).Value = -8246

# --- CRP row 89 (Leve Item ID 12584) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 17860142
$ws.Cells.Item(89, 10).Value = 6000
$ws.Cells.Item(89, 12).Value = 30000
$ws.Cells.Item(89, 14).Value = -41232

# --- CUL row 9 (Leve Item ID 4681) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 1254050.4
$ws.Cells.Item(9, 10).Value = 481.33334
$ws.Cells.Item(9, 12).Value = 1444.00002
$ws.Cells.Item(9, 14).Value = -1892.00002

# --- CUL row 54 (Leve Item ID 4671) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(54, 8).Value = 2990.1667
$ws.Cells.Item(54, 10).Value = 3533.2
$ws.Cells.Item(54, 12).Value = 10599.6
$ws.Cells.Item(54, 14).Value = -11717.6

# --- CUL row 134 (Leve Item ID 44074) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 1347.7646
$ws.Cells.Item(134, 9).Value = 991.93335
$ws.Cells.Item(134, 11).Value = 2975.80005
$ws.Cells.Item(134, 13).Value = 2094.19995

# --- CUL row 140 (Leve Item ID 44097) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 2133.5
$ws.Cells.Item(140, 9).Value = 1970.6923
$ws.Cells.Item(140, 11).Value = 5912.0769
$ws.Cells.Item(140, 13).Value = -732.0769

# --- GSM row 11 (Leve Item ID 4422) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 11155197
$ws.Cells.Item(11, 9).Value = 1119242.2
$ws.Cells.Item(11, 10).Value = 17845834
$ws.Cells.Item(11, 11).Value = 1119242.2
$ws.Cells.Item(11, 12).Value = 17845834
$ws.Cells.Item(11, 13).Value = -1119103.2
$ws.Cells.Item(11, 14).Value = -17846112

# --- GSM row 13 (Leve Item ID 2443) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 100
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 14).Value = ""

# --- GSM row 52 (Leve Item ID 4147) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 25200
$ws.Cells.Item(52, 10).Value = 26500
$ws.Cells.Item(52, 12).Value = 26500
$ws.Cells.Item(52, 14).Value = -27018

# --- GSM row 58 (Leve Item ID 4363) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(58, 8).Value = 26250
$ws.Cells.Item(58, 10).Value = 25000
$ws.Cells.Item(58, 12).Value = 25000
$ws.Cells.Item(58, 14).Value = -25554

# --- GSM row 93 (Leve Item ID 18107) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(93, 8).Value = 14167.8
$ws.Cells.Item(93, 10).Value = 14167.8
$ws.Cells.Item(93, 12).Value = 14167.8
$ws.Cells.Item(93, 14).Value = -17911.8

# --- GSM row 122 (Leve Item ID 36182) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 5105677.5
$ws.Cells.Item(122, 9).Value = 5910600
$ws.Cells.Item(122, 11).Value = 17731800
$ws.Cells.Item(122, 13).Value = -17729350

# --- GSM row 132 (Leve Item ID 44008) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3367.5676
$ws.Cells.Item(132, 9).Value = 2703.3333
$ws.Cells.Item(132, 10).Value = 6214.2856
$ws.Cells.Item(132, 11).Value = 8109.999899999999
$ws.Cells.Item(132, 12).Value = 18642.8568
$ws.Cells.Item(132, 13).Value = -5579.999899999999
$ws.Cells.Item(132, 14).Value = -23702.8568

# --- LTW row 46 (Leve Item ID 5282) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 5874.793
$ws.Cells.Item(46, 9).Value = 8239.532999999999
$ws.Cells.Item(46, 10).Value = 3341.1428
$ws.Cells.Item(46, 11).Value = 8239.532999999999
$ws.Cells.Item(46, 12).Value = 3341.1428
$ws.Cells.Item(46, 13).Value = -8051.532999999999
$ws.Cells.Item(46, 14).Value = -3717.1428

# --- LTW row 61 (Leve Item ID 27740) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1082
$ws.Cells.Item(61, 9).Value = 578.1667
$ws.Cells.Item(61, 10).Value = 2089.6667
$ws.Cells.Item(61, 11).Value = 578.1667
$ws.Cells.Item(61, 12).Value = 2089.6667
$ws.Cells.Item(61, 13).Value = -376.1667
$ws.Cells.Item(61, 14).Value = -2493.6667

# --- LTW row 82 (Leve Item ID 12565) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2520.2
$ws.Cells.Item(82, 9).Value = 2200.6667
$ws.Cells.Item(82, 10).Value = 2999.5
$ws.Cells.Item(82, 11).Value = 2200.6667
$ws.Cells.Item(82, 12).Value = 2999.5
$ws.Cells.Item(82, 13).Value = -1839.6667
$ws.Cells.Item(82, 14).Value = -3721.5

# --- LTW row 85 (Leve Item ID 12565) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 2520.2
$ws.Cells.Item(85, 9).Value = 2200.6667
$ws.Cells.Item(85, 10).Value = 2999.5
$ws.Cells.Item(85, 11).Value = 2200.6667
$ws.Cells.Item(85, 12).Value = 2999.5
$ws.Cells.Item(85, 13).Value = -952.6667000000002
$ws.Cells.Item(85, 14).Value = -5495.5

# --- LTW row 113 (Leve Item ID 27740) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 1082
$ws.Cells.Item(113, 9).Value = 578.1667
$ws.Cells.Item(113, 10).Value = 2089.6667
$ws.Cells.Item(113, 11).Value = 578.1667
$ws.Cells.Item(113, 12).Value = 2089.6667
$ws.Cells.Item(113, 13).Value = 1591.8333
$ws.Cells.Item(113, 14).Value = -6429.6667

# --- LTW row 132 (Leve Item ID 44058) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3765.8572
$ws.Cells.Item(132, 9).Value = 3384.6
$ws.Cells.Item(132, 10).Value = 3977.6667
$ws.Cells.Item(132, 11).Value = 10153.8
$ws.Cells.Item(132, 12).Value = 11933.0001
$ws.Cells.Item(132, 13).Value = -7623.799999999999
$ws.Cells.Item(132, 14).Value = -16993.0001

# --- LTW row 136 (Leve Item ID 44060) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2006.0526
$ws.Cells.Item(136, 9).Value = 1675.1333
$ws.Cells.Item(136, 10).Value = 3247
$ws.Cells.Item(136, 11).Value = 5025.3999
$ws.Cells.Item(136, 12).Value = 9741
$ws.Cells.Item(136, 13).Value = -2475.3999
$ws.Cells.Item(136, 14).Value = -14841

# --- LTW row 139 (Leve Item ID 43310) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).Value = ""

# --- WVR row 80 (Leve Item ID 10911) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(80, 8).Value = 75000
$ws.Cells.Item(80, 10).Value = 75000
$ws.Cells.Item(80, 12).Value = 75000
$ws.Cells.Item(80, 14).Value = -76996

# --- WVR row 83 (Leve Item ID 10911) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(83, 8).Value = 75000
$ws.Cells.Item(83, 10).Value = 75000
$ws.Cells.Item(83, 12).Value = 225000
$ws.Cells.Item(83, 14).Value = -234984

# --- WVR row 132 (Leve Item ID 44029) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2663.0557
$ws.Cells.Item(132, 9).Value = 2418.2307
$ws.Cells.Item(132, 10).Value = 3299.6
$ws.Cells.Item(132, 11).Value = 7254.6921
$ws.Cells.Item(132, 12).Value = 9898.799999999999
$ws.Cells.Item(132, 13).Value = -4724.6921
$ws.Cells.Item(132, 14).Value = -14958.8
